$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header "Amount" -> "Quantity"
$ws.Range("B1").Value = "Quantity"

# Update MSFT quantity 14 -> 5
$ws.Range("B2").Value = 5

# Match the final selection shown in the saved file (B3)
$ws.Range("B3").Select()
